$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that change from CHARTER to STAY
$toStay = @(27, 28, 82, 83, 84, 85, 86, 87)
foreach ($r in $toStay) {
    $ws.Range("A$r").Value = "STAY"
}

# Rows that change from STAY to CHARTER
$toCharter = @(110, 114, 115, 116)
foreach ($r in $toCharter) {
    $ws.Range("A$r").Value = "CHARTER"
}
